$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
    # Row 12
    $ws.Range("H12").Value = 127
    $ws.Range("I12").Value = 123.166664
    $ws.Range("K12").Value = 123.166664
    $ws.Range("M12").Value = 46.833336
    # Row 17
    $ws.Range("H17").Value = 2394759
    $ws.Range("J17").Value = 2431042
    $ws.Range("L17").Value = 7293126
    $ws.Range("N17").Value = -7293462
    # Row 64
    $ws.Range("H64").Value = 7900
    $ws.Range("I64").Value = 8371.429
    $ws.Range("J64").Value = 4600
    $ws.Range("K64").Value = 8371.429
    $ws.Range("L64").Value = 4600
    $ws.Range("M64").Value = -8123.429
    $ws.Range("N64").Value = -5096
    # Row 67
    $ws.Range("H67").Value = 7900
    $ws.Range("I67").Value = 8371.429
    $ws.Range("J67").Value = 4600
    $ws.Range("K67").Value = 8371.429
    $ws.Range("L67").Value = 4600
    $ws.Range("M67").Value = -7513.429
    $ws.Range("N67").Value = -6316
    # Row 69
    $ws.Range("H69").Value = 5220
    $ws.Range("I69").Value = 4800
    $ws.Range("J69").Value = 5340
    $ws.Range("K69").Value = 14400
    $ws.Range("L69").Value = 16020
    $ws.Range("M69").Value = -13526
    $ws.Range("N69").Value = -17768
    # Row 72
    $ws.Range("H72").Value = 5220
    $ws.Range("I72").Value = 4800
    $ws.Range("J72").Value = 5340
    $ws.Range("K72").Value = 43200
    $ws.Range("L72").Value = 48060
    $ws.Range("M72").Value = -38832
    $ws.Range("N72").Value = -56796
    # Row 86
    $ws.Range("H86").Value = 2700
    $ws.Range("I86").Value = 0
    $ws.Range("J86").Value = 2700
    $ws.Range("K86").Value = 0
    $ws.Range("L86").Value = 2700
    $ws.Range("M86").ClearContents()
    $ws.Range("N86").Value = -4946
    # Row 89
    $ws.Range("H89").Value = 2700
    $ws.Range("I89").Value = 0
    $ws.Range("J89").Value = 2700
    $ws.Range("K89").Value = 0
    $ws.Range("L89").Value = 13500
    $ws.Range("M89").ClearContents()
    $ws.Range("N89").Value = -24732
    # Row 138
    $ws.Range("H138").Value = 1969.885
    $ws.Range("I138").Value = 614.5714
    $ws.Range("J138").Value = 4418.1934
    $ws.Range("K138").Value = 1843.7142
    $ws.Range("L138").Value = 13254.5802
    $ws.Range("M138").Value = 3296.2858
    $ws.Range("N138").Value = -23534.5802

$ws = $wb.Worksheets.Item("ARM")
    # Row 74
    $ws.Range("H74").Value = 13515008
    $ws.Range("I74").Value = 1312.1333
    $ws.Range("J74").Value = 71430850
    $ws.Range("K74").Value = 1312.1333
    $ws.Range("L74").Value = 71430850
    $ws.Range("M74").Value = -438.1333
    $ws.Range("N74").Value = -71432598
    # Row 77
    $ws.Range("H77").Value = 13515008
    $ws.Range("I77").Value = 1312.1333
    $ws.Range("J77").Value = 71430850
    $ws.Range("K77").Value = 6560.666499999999
    $ws.Range("L77").Value = 357154250
    $ws.Range("M77").Value = -2192.666499999999
    $ws.Range("N77").Value = -357162986
    # Row 122
    $ws.Range("H122").Value = 1511719.8
    $ws.Range("I122").Value = 1835038.2
    $ws.Range("J122").Value = 2900
    $ws.Range("K122").Value = 5505114.6
    $ws.Range("L122").Value = 8700
    $ws.Range("M122").Value = -5502664.6
    $ws.Range("N122").Value = -13600

$ws = $wb.Worksheets.Item("BSM")
    # Row 22
    $ws.Range("H22").Value = 2000
    $ws.Range("I22").Value = 2000
    $ws.Range("J22").Value = 0
    $ws.Range("K22").Value = 2000
    $ws.Range("L22").Value = 0
    $ws.Range("M22").Value = -1827
    $ws.Range("N22").ClearContents()
    # Row 96
    $ws.Range("H96").Value = 10341.6
    $ws.Range("I96").Value = 10341.6
    $ws.Range("J96").Value = 0
    $ws.Range("K96").Value = 10341.6
    $ws.Range("L96").Value = 0
    $ws.Range("M96").Value = -7595.6
    $ws.Range("N96").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
    # Row 22
    $ws.Range("H22").Value = 294.18182
    $ws.Range("I22").Value = 361
    $ws.Range("J22").Value = 116
    $ws.Range("K22").Value = 361
    $ws.Range("L22").Value = 116
    $ws.Range("M22").Value = -11
    $ws.Range("N22").Value = -816
    # Row 31
    $ws.Range("H31").Value = 5804286.5
    $ws.Range("I31").Value = 1305.8302
    $ws.Range("J31").Value = 14347564
    $ws.Range("K31").Value = 1305.8302
    $ws.Range("L31").Value = 14347564
    $ws.Range("M31").Value = -1010.8302
    $ws.Range("N31").Value = -14348154
    # Row 34
    $ws.Range("H34").Value = 5804286.5
    $ws.Range("I34").Value = 1305.8302
    $ws.Range("J34").Value = 14347564
    $ws.Range("K34").Value = 1305.8302
    $ws.Range("L34").Value = 14347564
    $ws.Range("M34").Value = -1103.8302
    $ws.Range("N34").Value = -14347968
    # Row 58
    $ws.Range("H58").Value = 2689290.2
    $ws.Range("I58").Value = 4505227
    $ws.Range("J58").Value = 1703.68
    $ws.Range("K58").Value = 4505227
    $ws.Range("L58").Value = 1703.68
    $ws.Range("M58").Value = -4505024
    $ws.Range("N58").Value = -2109.68
    # Row 136
    $ws.Range("H136").Value = 2689290.2
    $ws.Range("I136").Value = 4505227
    $ws.Range("J136").Value = 1703.68
    $ws.Range("K136").Value = 13515681
    $ws.Range("L136").Value = 5111.04
    $ws.Range("M136").Value = -13513131
    $ws.Range("N136").Value = -10211.04
    # Row 140
    $ws.Range("H140").Value = 27469.9
    $ws.Range("J140").Value = 27469.9
    $ws.Range("L140").Value = 27469.9
    $ws.Range("N140").Value = -37829.9

$ws = $wb.Worksheets.Item("CUL")
    # Row 121
    $ws.Range("H121").Value = 1011.1852
    $ws.Range("J121").Value = 1098.3478
    $ws.Range("L121").Value = 3295.0434
    $ws.Range("N121").Value = -5915.0434
    # Row 133
    $ws.Range("H133").Value = 47482.44
    $ws.Range("I133").Value = 105206.1
    $ws.Range("K133").Value = 315618.3
    $ws.Range("M133").Value = -310558.3

$ws = $wb.Worksheets.Item("GSM")
    # Row 80
    $ws.Range("H80").Value = 2684.24
    $ws.Range("I80").Value = 2329.4119
    $ws.Range("K80").Value = 2329.4119
    $ws.Range("M80").Value = -1331.4119
    # Row 83
    $ws.Range("H83").Value = 2684.24
    $ws.Range("I83").Value = 2329.4119
    $ws.Range("K83").Value = 11647.0595
    $ws.Range("M83").Value = -6655.059499999999
    # Row 122
    $ws.Range("H122").Value = 79631450
    $ws.Range("I122").Value = 106483130
    $ws.Range("J122").Value = 12502250
    $ws.Range("K122").Value = 319449390
    $ws.Range("L122").Value = 37506750
    $ws.Range("M122").Value = -319446940
    $ws.Range("N122").Value = -37511650

$ws = $wb.Worksheets.Item("LTW")
    # Row 61
    $ws.Range("H61").Value = 1856.5
    $ws.Range("I61").Value = 1707.4286
    $ws.Range("J61").Value = 2900
    $ws.Range("K61").Value = 1707.4286
    $ws.Range("L61").Value = 2900
    $ws.Range("M61").Value = -1505.4286
    $ws.Range("N61").Value = -3304
    # Row 68
    $ws.Range("H68").Value = 83335100
    $ws.Range("I68").Value = 1925
    $ws.Range("J68").Value = 250001460
    $ws.Range("K68").Value = 1925
    $ws.Range("L68").Value = 250001460
    $ws.Range("M68").Value = -1176
    $ws.Range("N68").Value = -250002958
    # Row 71
    $ws.Range("H71").Value = 83335100
    $ws.Range("I71").Value = 1925
    $ws.Range("J71").Value = 250001460
    $ws.Range("K71").Value = 9625
    $ws.Range("L71").Value = 1250007300
    $ws.Range("M71").Value = -5881
    $ws.Range("N71").Value = -1250014788
    # Row 82
    $ws.Range("H82").Value = 35751.1
    $ws.Range("I82").Value = 1499.0588
    $ws.Range("J82").Value = 80542.234
    $ws.Range("K82").Value = 1499.0588
    $ws.Range("L82").Value = 80542.234
    $ws.Range("M82").Value = -1138.0588
    $ws.Range("N82").Value = -81264.234
    # Row 85
    $ws.Range("H85").Value = 35751.1
    $ws.Range("I85").Value = 1499.0588
    $ws.Range("J85").Value = 80542.234
    $ws.Range("K85").Value = 1499.0588
    $ws.Range("L85").Value = 80542.234
    $ws.Range("M85").Value = -251.0588
    $ws.Range("N85").Value = -83038.234
    # Row 113
    $ws.Range("H113").Value = 1856.5
    $ws.Range("I113").Value = 1707.4286
    $ws.Range("J113").Value = 2900
    $ws.Range("K113").Value = 1707.4286
    $ws.Range("L113").Value = 2900
    $ws.Range("M113").Value = 462.5714
    $ws.Range("N113").Value = -7240
    # Row 132
    $ws.Range("H132").Value = 9527262
    $ws.Range("I132").Value = 10887744
    $ws.Range("J132").Value = 3886.1428
    $ws.Range("K132").Value = 32663232
    $ws.Range("L132").Value = 11658.4284
    $ws.Range("M132").Value = -32660702
    $ws.Range("N132").Value = -16718.4284

$ws = $wb.Worksheets.Item("WVR")
    # Row 62
    $ws.Range("H62").Value = 3150
    $ws.Range("I62").Value = 3150
    $ws.Range("K62").Value = 3150
    $ws.Range("M62").Value = -2526
    # Row 65
    $ws.Range("H65").Value = 3150
    $ws.Range("I65").Value = 3150
    $ws.Range("K65").Value = 15750
    $ws.Range("M65").Value = -12630
    # Row 132
    $ws.Range("H132").Value = 921.8280999999999
    $ws.Range("I132").Value = 520.89795
    $ws.Range("J132").Value = 2231.5334
    $ws.Range("K132").Value = 1562.69385
    $ws.Range("L132").Value = 6694.600199999999
    $ws.Range("M132").Value = 967.3061499999999
    $ws.Range("N132").Value = -11754.6002
    # Row 136
    $ws.Range("H136").Value = 12347266
    $ws.Range("I136").Value = 1421.575
    $ws.Range("J136").Value = 47621104
    $ws.Range("K136").Value = 4264.725
    $ws.Range("L136").Value = 142863312
    $ws.Range("M136").Value = -1714.725
    $ws.Range("N136").Value = -142868412
